$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 12 (existing rows 12-17 shift down to 14-19)
$ws.Range("A12:A13").EntireRow.Insert()

# --- New row 12: Alcachofa Argentina(o) ---
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C12").Value = "Los Lagos"
$ws.Range("D12").Value = 44757
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112013
$ws.Range("G12").Value = "Alcachofa"
$ws.Range("H12").Value = "Argentina(o)"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 18000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 18000
$ws.Range("N12").Value = "$/caja 50 unidades"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 360
$ws.Range("Q12").Value = 50
$ws.Range("R12").Value = "Hortaliza"

# --- New row 13: Alcachofa Española ---
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("D13").Value = 44757
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 100112013
$ws.Range("G13").Value = "Alcachofa"
$ws.Range("H13").Value = "Española"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 70
$ws.Range("K13").Value = 22000
$ws.Range("L13").Value = 22000
$ws.Range("M13").Value = 22000
$ws.Range("N13").Value = "$/caja 30 unidades"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 733
$ws.Range("Q13").Value = 30
$ws.Range("R13").Value = "Hortaliza"
